$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 updates
$ws.Range("G6").Value = 6.25
$ws.Range("H6").Value = 3.7
$ws.Range("I6").Value = 1.53
$ws.Range("J6").Value = 6.5
$ws.Range("L6").Value = 2.2
$ws.Range("U6").Value = 2.2
$ws.Range("V6").Value = 1.62
$ws.Range("Y6").Value = 21
$ws.Range("AD6").Value = 7.5
$ws.Range("AE6").Value = 21
$ws.Range("AF6").Value = 81
$ws.Range("AH6").Value = 5.5
$ws.Range("AI6").Value = 6.5
$ws.Range("AK6").Value = 11
$ws.Range("AN6").Value = 7.5
$ws.Range("AO6").Value = 34
$ws.Range("AQ6").Value = 151
$ws.Range("AU6").Value = 9.5
$ws.Range("AW6").Value = 3.4
$ws.Range("AX6").Value = 8
$ws.Range("AZ6").Value = 26

# Row 7 updates
$ws.Range("G7").Value = 1.62
